$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 179 (shifts existing rows 179-210 down to 180-211)
$ws.Rows.Item(179).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(179, 1).Value = 9
$ws.Cells.Item(179, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(179, 3).Value = "Metropolitana"
$ws.Cells.Item(179, 4).Value = 44694
$ws.Cells.Item(179, 5).Value = 13
$ws.Cells.Item(179, 6).Value = 100112026
$ws.Cells.Item(179, 7).Value = "Haba"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 61
$ws.Cells.Item(179, 11).Value = 16000
$ws.Cells.Item(179, 12).Value = 17000
$ws.Cells.Item(179, 13).Value = 16508
$ws.Cells.Item(179, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(179, 15).Value = "Provincia de Talagante"
$ws.Cells.Item(179, 16).Value = 660
$ws.Cells.Item(179, 17).Value = 25
$ws.Cells.Item(179, 18).Value = "Hortaliza"
